# MNX Track update: refresh the test-run tracking numbers on both sheets.
#
# "MXTMS" (sheet 1), row 2 is the live QA/Selenium smoke-test record:
#   D2 = OrderID, E2 = PickupID
# "Connect" (sheet 2), row 2 is the matching Connect-side record:
#   AG2 = PickupID
#
# All three cells store their (numeric-looking) identifiers as TEXT, not
# numbers, in the original workbook (<c t="s">). A plain `.Value =` on a
# digit-only string gets auto-coerced to a number by Excel, so we force
# the Text number format on each cell *before* writing the new value -
# this keeps the cells as shared-string/text cells exactly like before.

$wb = $excel.ActiveWorkbook

$wsConnect = $wb.Worksheets.Item("Connect")
$wsMxtms = $wb.Worksheets.Item("MXTMS")

# Connect!AG2 (PickupID) : 10016605 -> 15592291
$wsConnect.Range("AG2").NumberFormat = "@"
$wsConnect.Range("AG2").Value = "15592291"

# MXTMS!D2 (OrderID) : 11181635 -> 11189676
$wsMxtms.Range("D2").NumberFormat = "@"
$wsMxtms.Range("D2").Value = "11189676"

# MXTMS!E2 (PickupID) : 7386723 -> 7391426
$wsMxtms.Range("E2").NumberFormat = "@"
$wsMxtms.Range("E2").Value = "7391426"
